$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value  = "3DS Digi 3"
$ws.Range("A5").Value  = "3DS Digi 4"
$ws.Range("A6").Value  = "3DS Digi 5"
$ws.Range("A7").Value  = "3DS Digi 6"
$ws.Range("A8").Value  = "3DS Digi 7"

$ws.Range("A17").Value = "Switch Digi 1"
$ws.Range("A18").Value = "Switch Digi 2"
$ws.Range("A19").Value = "Switch Digi 3"

$ws.Range("A26").Value = "Wii U Digi 5"
$ws.Range("A27").Value = "Wii U Digi 6"
$ws.Range("A28").Value = "Wii U Digi 7"
$ws.Range("A29").Value = "Wii U Digi 8"
